$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Last Updated timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 10:20 AM"

# --- Top Gainers sheet ---
$wsGain = $wb.Worksheets.Item("Top Gainers")
$wsGain.Range("B17").Value = "HCG"
$wsGain.Range("C17").Value = 6.5614
$wsGain.Range("D17").Value = 2.2658
$wsGain.Range("E17").Value = 20.6746
$wsGain.Range("B18").Value = "POKARNA"
$wsGain.Range("C18").Value = 6.3931
$wsGain.Range("D18").Value = -1.0672
$wsGain.Range("E18").Value = 19.2413
$wsGain.Range("B19").Value = "SAIL"
$wsGain.Range("C19").Value = 6.3484
$wsGain.Range("D19").Value = 8.5663
$wsGain.Range("E19").Value = 4.5137
$wsGain.Range("B20").Value = "ABDL"
$wsGain.Range("C20").Value = 6.0222
$wsGain.Range("D20").Value = 4.9163
$wsGain.Range("E20").Value = 27.7822
$wsGain.Range("B21").Value = "WALCHANNAG"
$wsGain.Range("C21").Value = 6.0056
$wsGain.Range("D21").Value = 3.5575
$wsGain.Range("E21").Value = -5.588
$wsGain.Range("B22").Value = "TMB"
$wsGain.Range("C22").Value = 6.0054
$wsGain.Range("D22").Value = 9.7967
$wsGain.Range("E22").Value = 17.1347
$wsGain.Range("B23").Value = "RPOWER"
$wsGain.Range("C23").Value = 5.9982
$wsGain.Range("D23").Value = 2.4336
$wsGain.Range("E23").Value = 4.3733
$wsGain.Range("B24").Value = "CELLO"
$wsGain.Range("C24").Value = 5.9078
$wsGain.Range("D24").Value = 4.7408
$wsGain.Range("E24").Value = 14.6875
$wsGain.Range("B25").Value = "M&MFIN"
$wsGain.Range("C25").Value = 5.852
$wsGain.Range("D25").Value = 6.4733
$wsGain.Range("E25").Value = 15.3106
$wsGain.Range("B26").Value = "FISCHER"
$wsGain.Range("C26").Value = 5.7481
$wsGain.Range("D26").Value = 10.7325
$wsGain.Range("E26").Value = 3.9189
$wsGain.Range("B27").Value = "SRM"
$wsGain.Range("C27").Value = 5.7289
$wsGain.Range("D27").Value = 5.4147
$wsGain.Range("E27").Value = 6.3431
$wsGain.Range("B28").Value = "JISLJALEQS"
$wsGain.Range("C28").Value = 5.5917
$wsGain.Range("D28").Value = 4.8879
$wsGain.Range("E28").Value = -1.1798
$wsGain.Range("B29").Value = "STLTECH"
$wsGain.Range("C29").Value = 5.5556
$wsGain.Range("D29").Value = 2.5231
$wsGain.Range("E29").Value = 8.696
$wsGain.Range("B30").Value = "IOC"
$wsGain.Range("C30").Value = 5.5462
$wsGain.Range("D30").Value = 8.4591
$wsGain.Range("E30").Value = 8.8791
$wsGain.Range("B31").Value = "PDSL"
$wsGain.Range("C31").Value = 5.1546
$wsGain.Range("D31").Value = 3.1682
$wsGain.Range("E31").Value = 9.0052
$wsGain.Range("B32").Value = "ADANIENSOL"
$wsGain.Range("C32").Value = 5.0315
$wsGain.Range("D32").Value = 2.4676
$wsGain.Range("E32").Value = 10.9258
$wsGain.Range("B33").Value = "MEGASOFT"
$wsGain.Range("C33").Value = 4.9974
$wsGain.Range("D33").Value = 15.7588
$wsGain.Range("E33").Value = 33.5271
$wsGain.Range("B34").Value = "INDIANHUME"
$wsGain.Range("C34").Value = 4.9964
$wsGain.Range("D34").Value = 5.3443
$wsGain.Range("E34").Value = 1.8807
$wsGain.Range("B35").Value = "AXISCADES"
$wsGain.Range("C35").Value = 4.9963
$wsGain.Range("D35").Value = 7.474
$wsGain.Range("E35").Value = -2.5721
$wsGain.Range("B36").Value = "PROZONER"
$wsGain.Range("C36").Value = 4.9921
$wsGain.Range("D36").Value = 15.7468
$wsGain.Range("E36").Value = 36.095
$wsGain.Range("B37").Value = "STALLION"
$wsGain.Range("C37").Value = 4.9914
$wsGain.Range("D37").Value = -5.2229
$wsGain.Range("E37").Value = 21.4391
$wsGain.Range("B38").Value = "SURYAROSNI"
$wsGain.Range("C38").Value = 4.9831
$wsGain.Range("D38").Value = 11.386
$wsGain.Range("E38").Value = 3.0213
$wsGain.Range("B39").Value = "BIL"
$wsGain.Range("C39").Value = 4.9242
$wsGain.Range("D39").Value = 9.7065
$wsGain.Range("E39").Value = 0.2135
$wsGain.Range("B40").Value = "DATAMATICS"
$wsGain.Range("C40").Value = 4.9005
$wsGain.Range("D40").Value = 7.3152
$wsGain.Range("E40").Value = 15.7298
$wsGain.Range("B41").Value = "UTKARSHBNK"
$wsGain.Range("C41").Value = 4.8768
$wsGain.Range("D41").Value = -5.8959
$wsGain.Range("E41").Value = -2.6215
$wsGain.Range("B42").Value = "FILATEX"
$wsGain.Range("C42").Value = 4.8689
$wsGain.Range("D42").Value = 10.274
$wsGain.Range("E42").Value = 26.0027
$wsGain.Range("B43").Value = "HITECHGEAR"
$wsGain.Range("C43").Value = 4.8651
$wsGain.Range("D43").Value = 2.1287
$wsGain.Range("E43").Value = 10.9905
$wsGain.Range("B44").Value = "INDOTHAI"
$wsGain.Range("C44").Value = 4.8064
$wsGain.Range("D44").Value = 4.5349
$wsGain.Range("E44").Value = 43.748
$wsGain.Range("B45").Value = "SAPPHIRE"
$wsGain.Range("C45").Value = 4.6445
$wsGain.Range("D45").Value = 2.2696
$wsGain.Range("E45").Value = -0.3063
$wsGain.Range("B74").Value = "BCLIND"
$wsGain.Range("C74").Value = 3.6271
$wsGain.Range("D74").Value = 2.2945
$wsGain.Range("E74").Value = 0.1728
$wsGain.Range("B75").Value = "CGPOWER"
$wsGain.Range("C75").Value = 3.6125
$wsGain.Range("D75").Value = 3.4192
$wsGain.Range("E75").Value = 1.0325
$wsGain.Range("B76").Value = "WELSPUNLIV"
$wsGain.Range("C76").Value = 3.6073
$wsGain.Range("D76").Value = 3.7285
$wsGain.Range("E76").Value = 15.9372

# --- Top Losers sheet ---
$wsLose = $wb.Worksheets.Item("Top Losers")
$wsLose.Range("B11").Value = "TVSELECT"
$wsLose.Range("C11").Value = -5.6153
$wsLose.Range("D11").Value = -0.9738
$wsLose.Range("E11").Value = -2.9968
$wsLose.Range("B12").Value = "NSLNISP"
$wsLose.Range("C12").Value = -5.4542
$wsLose.Range("D12").Value = 1.3037
$wsLose.Range("E12").Value = 0.4681
$wsLose.Range("B30").Value = "BHARATWIRE"
$wsLose.Range("C30").Value = -3.5327
$wsLose.Range("D30").Value = 22.8336
$wsLose.Range("E30").Value = 23.8979
$wsLose.Range("B31").Value = "ABSLAMC"
$wsLose.Range("C31").Value = -3.5313
$wsLose.Range("D31").Value = -5.9355
$wsLose.Range("E31").Value = -1.2887
$wsLose.Range("B32").Value = "SPLPETRO"
$wsLose.Range("C32").Value = -3.3984
$wsLose.Range("D32").Value = -5.0241
$wsLose.Range("E32").Value = -7.769
$wsLose.Range("B33").Value = "CAMS"
$wsLose.Range("C33").Value = -3.2545
$wsLose.Range("D33").Value = -0.6366
$wsLose.Range("E33").Value = 2.5781
$wsLose.Range("B34").Value = "PRUDENT"
$wsLose.Range("C34").Value = -3.2484
$wsLose.Range("D34").Value = -3.6312
$wsLose.Range("E34").Value = 1.9933
$wsLose.Range("B35").Value = "SPARC"
$wsLose.Range("C35").Value = -3.1709
$wsLose.Range("D35").Value = 4.8337
$wsLose.Range("E35").Value = 6.3311
$wsLose.Range("B36").Value = "ANANDRATHI"
$wsLose.Range("C36").Value = -3.0775
$wsLose.Range("D36").Value = -0.8672
$wsLose.Range("E36").Value = 9.1835
$wsLose.Range("B37").Value = "NLCINDIA"
$wsLose.Range("C37").Value = -3.0757
$wsLose.Range("D37").Value = -4.5618
$wsLose.Range("E37").Value = -11.6431
$wsLose.Range("B38").Value = "YATRA"
$wsLose.Range("C38").Value = -3.0403
$wsLose.Range("D38").Value = -2.8455
$wsLose.Range("E38").Value = 7.3711
$wsLose.Range("B39").Value = "MPSLTD"
$wsLose.Range("C39").Value = -3.0335
$wsLose.Range("D39").Value = -4.3902
$wsLose.Range("E39").Value = 2.434
$wsLose.Range("B40").Value = "DRREDDY"
$wsLose.Range("C40").Value = -2.9859
$wsLose.Range("D40").Value = -2.5475
$wsLose.Range("E40").Value = 2.2228
$wsLose.Range("B41").Value = "ROSSTECH"
$wsLose.Range("C41").Value = -2.9778
$wsLose.Range("D41").Value = 1.9028
$wsLose.Range("E41").Value = -6.8057
$wsLose.Range("B42").Value = "OAL"
$wsLose.Range("C42").Value = -2.9496
$wsLose.Range("D42").Value = -1.278
$wsLose.Range("E42").Value = 8.7362
$wsLose.Range("B43").Value = "ENDURANCE"
$wsLose.Range("C43").Value = -2.939
$wsLose.Range("D43").Value = -2.2945
$wsLose.Range("E43").Value = 3.4531
$wsLose.Range("B44").Value = "POLICYBZR"
$wsLose.Range("C44").Value = -2.907
$wsLose.Range("D44").Value = 2.2365
$wsLose.Range("E44").Value = 1.2573
$wsLose.Range("B45").Value = "BOSCHLTD"
$wsLose.Range("C45").Value = -2.9061
$wsLose.Range("D45").Value = -3.0193
$wsLose.Range("E45").Value = -1.9006
$wsLose.Range("B46").Value = "AYMSYNTEX"
$wsLose.Range("C46").Value = -2.9052
$wsLose.Range("D46").Value = -0.3705
$wsLose.Range("E46").Value = -10.494

# --- 1 Month Performance sheet ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Range("B21").Value = "SAMMAANCAP"
$wsPerf.Range("C21").Value = 35.5128
$wsPerf.Range("B22").Value = "SOUTHBANK"
$wsPerf.Range("C22").Value = 35.2819
$wsPerf.Range("B23").Value = "SHAREINDIA"
$wsPerf.Range("C23").Value = 35.2728
$wsPerf.Range("B24").Value = "TVSELECT"
$wsPerf.Range("C24").Value = 35.1983
$wsPerf.Range("B40").Value = "HATSUN"
$wsPerf.Range("C40").Value = 26.492
$wsPerf.Range("B41").Value = "ADANIPOWER"
$wsPerf.Range("C41").Value = 25.8247
